# My entire MTG collection is worth around 1230 EUR
# Refresh the price-tracker sheet with the latest pull of the collection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old (smaller) card list in A2:D12 with the new, larger one in A2:D24.
$ws.Range("A2").Value = "Ojer Taq, Deepest Foundation // Temple of Civilization"
$ws.Range("B2").Value = "The Lost Caverns of Ixalan: Extras"
$ws.Range("C2").Value = "Foil"
$ws.Range("D2").Value = 12.3

$ws.Range("A3").Value = "Blazing Archon"
$ws.Range("B3").Value = "Ravnica Remastered: Extras"
$ws.Range("C3").Value = "Normal"
$ws.Range("D3").Value = 1.41

$ws.Range("A4").Value = "Ojer Pakpatiq, Deepest Epoch // Temple of Cyclical Time"
$ws.Range("B4").Value = "The Lost Caverns of Ixalan: Extras"
$ws.Range("C4").Value = "Normal"
$ws.Range("D4").Value = 1.92

$ws.Range("A5").Value = "Cyclonic Rift"
$ws.Range("B5").Value = "Ravnica Remastered: Extras"
$ws.Range("C5").Value = "V.2"
$ws.Range("D5").Value = 38.53

$ws.Range("A6").Value = "Profane Tutor"
$ws.Range("B6").Value = "Modern Horizons 2"
$ws.Range("C6").Value = "Normal"
$ws.Range("D6").Value = 2.39

$ws.Range("A7").Value = "Stalactite Stalker"
$ws.Range("B7").Value = "The Lost Caverns of Ixalan"
$ws.Range("C7").Value = "Normal"
$ws.Range("D7").Value = 1.07

$ws.Range("A8").Value = "Trumpeting Carnasaur"
$ws.Range("B8").Value = "The Lost Caverns of Ixalan"
$ws.Range("C8").Value = "Normal"
$ws.Range("D8").Value = 2.53

$ws.Range("A9").Value = "Hardened Scales"
$ws.Range("B9").Value = "Enchanting Tales"
$ws.Range("C9").Value = "Normal"
$ws.Range("D9").Value = 1.59

$ws.Range("A10").Value = "Likeness Looter"
$ws.Range("B10").Value = "Wilds of Eldraine: Promos"
$ws.Range("C10").Value = "V.1"
$ws.Range("D10").Value = 1.64

$ws.Range("A11").Value = "Eriette of the Charmed Apple"
$ws.Range("B11").Value = "Wilds of Eldraine"
$ws.Range("C11").Value = "Normal"
$ws.Range("D11").Value = 2.52

$ws.Range("A12").Value = "Kellan, the Fae-Blooded // Birthright Boon"
$ws.Range("B12").Value = "Wilds of Eldraine"
$ws.Range("C12").Value = "Normal"
$ws.Range("D12").Value = 2.2

$ws.Range("A13").Value = "Questing Druid"
$ws.Range("B13").Value = "Wilds of Eldraine: Promos"
$ws.Range("C13").Value = "V.2"
$ws.Range("D13").Value = 5.68

$ws.Range("A14").Value = "Molten Collapse"
$ws.Range("B14").Value = "The Lost Caverns of Ixalan"
$ws.Range("C14").Value = "Normal"
$ws.Range("D14").Value = 2.08

$ws.Range("A15").Value = "Molten Collapse"
$ws.Range("B15").Value = "The Lost Caverns of Ixalan: Extras"
$ws.Range("C15").Value = "Normal"
$ws.Range("D15").Value = 1.63

$ws.Range("A16").Value = "Pantlaza, Sun-Favored"
$ws.Range("B16").Value = "Commander: The Lost Caverns of Ixalan: Extras"
$ws.Range("C16").Value = "V.1 Foil"
$ws.Range("D16").Value = 7.7

$ws.Range("A17").Value = "Niv-Mizzet, the Firemind"
$ws.Range("B17").Value = "Release Promos"
$ws.Range("C17").Value = "Normal"
$ws.Range("D17").Value = 2

$ws.Range("A18").Value = "Mayhem Devil"
$ws.Range("B18").Value = "Ravnica Remastered"
$ws.Range("C18").Value = "Normal"
$ws.Range("D18").Value = 1.7

$ws.Range("A19").Value = "Lavinia, Azorius Renegade"
$ws.Range("B19").Value = "Ravnica Remastered: Extras"
$ws.Range("C19").Value = "Foil"
$ws.Range("D19").Value = 4.68

$ws.Range("A20").Value = "Scion of Draco"
$ws.Range("B20").Value = "Modern Horizons 2"
$ws.Range("C20").Value = "Normal"
$ws.Range("D20").Value = 2.19

$ws.Range("A21").Value = "Roaming Throne"
$ws.Range("B21").Value = "The Lost Caverns of Ixalan: Promos"
$ws.Range("C21").Value = "V.1"
$ws.Range("D21").Value = 13.07

$ws.Range("A22").Value = "Urza's Incubator"
$ws.Range("B22").Value = "Dominaria Remastered: Extras"
$ws.Range("C22").Value = "V.2"
$ws.Range("D22").Value = 20.64

$ws.Range("A23").Value = "Cloudstone Curio"
$ws.Range("B23").Value = "Ravnica Remastered: Extras"
$ws.Range("C23").Value = "V.2 Foil"
$ws.Range("D23").Value = 25.42

$ws.Range("A24").Value = "Plaza of Heroes"
$ws.Range("B24").Value = "Dominaria United"
$ws.Range("C24").Value = "Normal"
$ws.Range("D24").Value = 6.96

# The running total moves down to row 26 and now sums the wider D2:D24 range.
$ws.Range("D26").Formula = "=SUM(D2:D24)"

# Re-fit the Name/Set columns now that several entries are much longer.
$ws.Columns.Item(1).ColumnWidth = 51.166666666666664
$ws.Columns.Item(2).ColumnWidth = 42.166666666666664

# Leave the selection on the new total cell.
[void]$ws.Range("D25").Select()
